$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.050.47"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "2.915.95"
$ws.Range("E3").Value = "  -0.44%  "
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "587.66"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -1.18%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "146.36"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +3.28%  "
$ws.Range("E7").Value = "  +0.13%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.504"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +0.79%  "
$ws.Range("D9").Value = "2.914.90"
$ws.Range("E9").Value = "  -0.49%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "7.06"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -1.47%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.151"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +6.80%  "
$ws.Range("E12").Value = "  -1.61%  "
$ws.Range("E13").Value = "  +6.91%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "32.32"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.76%  "
$ws.Range("E15").Value = "  -1.26%  "
$ws.Range("D16").Value = "3.400.70"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").Value = "62.048.39"
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").Value = "2.916.43"
$ws.Range("E19").Value = "  -1.82%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "432.86"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("E21").Value = "  -0.63%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.658"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -1.36%  "
$ws.Range("E23").Value = "  -1.82%  "
$ws.Range("B24").Value = "RenderToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "11.05"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +5.05%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "79.96"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.99%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "11.88"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +2.15%  "
$ws.Range("E27").Value = "  -0.75%  "
$ws.Range("E28").Value = "  -0.15%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "7.28"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +6.65%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.0000103"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +20.19%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "2.57"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.29%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "2.11"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +0.92%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.109"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +3.25%  "
$ws.Range("E34").Value = "  +0.19%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "25.93"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -0.56%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.977"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -0.50%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "3.12"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +9.92%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "5.51"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -0.46%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "49.19"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +0.01%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "1.99"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +2.77%  "
$ws.Range("E41").Value = "  -1.37%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.274"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +0.95%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "39.19"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +1.25%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "135.29"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +2.05%  "
$ws.Range("D46").Value = "2.692.08"
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("E47").Value = "  +0.48%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "349.78"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -2.36%  "
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("E50").Value = "  +0.34%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "22.46"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -0.59%  "
